$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the forecast dates in column A (rows 3-13) from consecutive days
# to the first day of each subsequent month.
$ws.Range("A3").Value = (Get-Date -Year 2025 -Month 2 -Day 1).Date
$ws.Range("A4").Value = (Get-Date -Year 2025 -Month 3 -Day 1).Date
$ws.Range("A5").Value = (Get-Date -Year 2025 -Month 4 -Day 1).Date
$ws.Range("A6").Value = (Get-Date -Year 2025 -Month 5 -Day 1).Date
$ws.Range("A7").Value = (Get-Date -Year 2025 -Month 6 -Day 1).Date
$ws.Range("A8").Value = (Get-Date -Year 2025 -Month 7 -Day 1).Date
$ws.Range("A9").Value = (Get-Date -Year 2025 -Month 8 -Day 1).Date
$ws.Range("A10").Value = (Get-Date -Year 2025 -Month 9 -Day 1).Date
$ws.Range("A11").Value = (Get-Date -Year 2025 -Month 10 -Day 1).Date
$ws.Range("A12").Value = (Get-Date -Year 2025 -Month 11 -Day 1).Date
$ws.Range("A13").Value = (Get-Date -Year 2025 -Month 12 -Day 1).Date

# Column A width change (slightly wider, matching the resaved workbook)
$ws.Columns.Item(1).ColumnWidth = 27.8

# Re-apply the "no fill" interior on the header style (L1) so Excel drops the
# stale applyFill flag on that cell's style record.
$ws.Range("L1").Interior.ColorIndex = -4142
$ws.Range("L1").Interior.Pattern = -4142

# Selection moved to L22 (matches the last active cell noted in the file)
$ws.Range("L22").Select()
